$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.224.70'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('E2').Style = "Normal"
# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.690.11'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E3').Style = "Normal"
# Row 4
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('E4').Style = "Normal"
# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.79'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('E5').Style = "Normal"
# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.521'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.48%  '
$ws.Range('E6').Style = "Normal"
# Row 7
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E7').Style = "Normal"
# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '23.15'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +10.18%  '
$ws.Range('E8').Style = "Normal"
# Row 9
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +4.63%  '
$ws.Range('E9').Style = "Normal"
# Row 10
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.30%  '
$ws.Range('E10').Style = "Normal"
# Row 11
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('E11').Style = "Normal"
# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.927.39'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('E12').Style = "Normal"
# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.687.92'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.12%  '
$ws.Range('E13').Style = "Normal"
# Row 14
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +2.55%  '
$ws.Range('E14').Style = "Normal"
# Row 15
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +4.31%  '
$ws.Range('E15').Style = "Normal"
# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '67.20'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.95%  '
$ws.Range('E16').Style = "Normal"
# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '27.219.94'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('E17').Style = "Normal"
# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '236.89'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('E18').Style = "Normal"
# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '8.10'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -1.17%  '
$ws.Range('E19').Style = "Normal"
# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0746'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('E20').Style = "Normal"
# Row 21
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E21').Style = "Normal"
# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.58'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +3.04%  '
$ws.Range('E22').Style = "Normal"
# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.63'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +4.38%  '
$ws.Range('E23').Style = "Normal"
# Row 24
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -1.39%  '
$ws.Range('E24').Style = "Normal"
# Row 25
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('E25').Style = "Normal"
# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.35'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +1.55%  '
$ws.Range('E26').Style = "Normal"
# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.47'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +2.43%  '
$ws.Range('E27').Style = "Normal"
# Row 28
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.42%  '
$ws.Range('E28').Style = "Normal"
# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.23%  '
$ws.Range('E29').Style = "Normal"
# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0506'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.93%  '
$ws.Range('E30').Style = "Normal"
# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.18'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('E31').Style = "Normal"
# Row 32
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +1.73%  '
$ws.Range('E32').Style = "Normal"
# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.552.59'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +3.63%  '
$ws.Range('E33').Style = "Normal"
# Row 34
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +2.12%  '
$ws.Range('E34').Style = "Normal"
# Row 35
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('E35').Style = "Normal"
# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.606'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +2.98%  '
$ws.Range('E36').Style = "Normal"
# Row 37
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +3.57%  '
$ws.Range('E37').Style = "Normal"
# Row 39
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.32%  '
$ws.Range('E39').Style = "Normal"
# Row 40
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +1.65%  '
$ws.Range('E40').Style = "Normal"
# Row 41
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.78'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.37%  '
$ws.Range('E41').Style = "Normal"
# Row 42
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '69.41'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +2.57%  '
$ws.Range('E42').Style = "Normal"
# Row 43
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('E43').Style = "Normal"
# Row 44
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -0.65%  '
$ws.Range('E44').Style = "Normal"
# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.835.73'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('E45').Style = "Normal"
# Row 46
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +1.43%  '
$ws.Range('E46').Style = "Normal"
# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '90.76'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('E47').Style = "Normal"
# Row 48
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0₆0111'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +7.30%  '
$ws.Range('E48').Style = "Normal"
# Row 49
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.63'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +6.71%  '
$ws.Range('E49').Style = "Normal"
# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.34'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +7.26%  '
$ws.Range('E50').Style = "Normal"
# Row 51
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.16%  '
$ws.Range('E51').Style = "Normal"

Write-Output "Applied 85 cell updates"
